# "feat: separando em backend e frontend"
# The "produto" column (D: produto / Produto A..T) is dropped from the
# sheet entirely - it's being split out into a separate backend concern,
# leaving only email, data, valor, quantidade, categoria on the sheet.
# Deleting the whole column shifts quantidade (old E) into D and
# categoria (old F) into E, and Excel automatically prunes the now-unused
# "produto"/"Produto A".."Produto T" shared strings on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "produto" column (D) - cells to its right shift left.
$ws.Columns.Item(4).Delete()

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("F7").Select() | Out-Null
